# Updated cryptos list (price / volume refresh), mirroring the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: several Price (column D) values are plain decimal-looking text
# (e.g. "1.00", "7.30") that must stay as literal text (matching the
# original inlineStr cells) rather than being auto-coerced to numbers
# (which would silently drop trailing zeros, e.g. "1.00" -> 1). For
# those cells we briefly force a Text number format before writing the
# value, then restore the "Normal" style so no stray formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.186.46"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "2.365.90"
$ws.Range("E3").Value = "  +6.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "

$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("E13").Value = "  +13.36%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.11%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.105"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "2.717.28"
$ws.Range("E16").Value = "  +6.27%  "

$ws.Range("D17").Value = "2.522.13"
$ws.Range("E17").Value = "  +12.26%  "

$ws.Range("D18").Value = "43.171.32"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "253.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.85%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.57%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("E32").Value = "  -2.62%  "

$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.43%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.58%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.31%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0376"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.81%  "

$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("E43").Value = "  -2.24%  "

$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.03%  "

$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0994"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.497.54"
$ws.Range("E51").Value = "  +4.41%  "
